# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets
#
# 1. Insert a new "Player Info" sheet in front of the existing sheets.
# 2. Rename the MATCH_CARD_LINK column to MATCH_CODE on both the
#    "ODI Batting" and "ODI Bowling" sheets, and replace the full
#    scorecard URL values with just the trailing numeric match code.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Player Info" sheet before "ODI Batting" -----------
$battingSheetRef = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheetRef)
$playerInfo.Name = "Player Info"

# Re-fetch the sheet references by name now that the insertion shifted
# worksheet positions (an Add() before a sheet can re-seat prior object
# references positionally, so grab fresh handles by name here).
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Cells holding a purely-numeric match code must stay text (matching the
# source data's inlineStr cells), so force text entry via NumberFormat
# "@" and then drop the leftover number-format styling with ClearFormats
# so the cell ends up with no special style, just like the other plain
# data cells on these sheets.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $playerInfo.Range("A2") "6036"
$playerInfo.Range("B2").Value = "Brydon Alexander Carse"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast"

# --- 2. "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE --------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{
    2  = "4472"
    3  = "4473"
    4  = "4476"
    5  = "4599"
    6  = "4602"
    7  = "4609"
    8  = "4613"
    9  = "4618"
    10 = "4619"
}
foreach ($row in $battingCodes.Keys) {
    Set-TextValue $battingSheet.Range("D$row") $battingCodes[$row]
}

# --- 3. "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE --------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{
    2  = "4472"
    3  = "4473"
    4  = "4476"
    5  = "4599"
    6  = "4602"
    7  = "4609"
    8  = "4613"
    9  = "4618"
    10 = "4619"
}
foreach ($row in $bowlingCodes.Keys) {
    Set-TextValue $bowlingSheet.Range("B$row") $bowlingCodes[$row]
}
